# Apply the commit's changes to the active worksheet:
#  1. Add a "pair_kind" value ("generic") to column J for rows 2-5.
#  2. Append a new "stim details" block starting at row 27:
#     - A27: "stim details"
#     - Row 28: header row (month, word_type, need_audio, need_image, word, count, find images)
#     - Rows 29-32: count=6,6,7,7 / type="video"
#     - Rows 33-36: count=6,6,7,7 / type="audio"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New "pair_kind" (column J) values for the practice-pair rows ---
$ws.Range("J2").Value = "generic"
$ws.Range("J3").Value = "generic"
$ws.Range("J4").Value = "generic"
$ws.Range("J5").Value = "generic"

# --- 2. New "stim details" block at the bottom of the sheet ---
$ws.Range("A27").Value = "stim details"

$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

$ws.Range("A29").Value = 6
$ws.Range("B29").Value = "video"

$ws.Range("A30").Value = 6
$ws.Range("B30").Value = "video"

$ws.Range("A31").Value = 7
$ws.Range("B31").Value = "video"

$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "video"

$ws.Range("A33").Value = 6
$ws.Range("B33").Value = "audio"

$ws.Range("A34").Value = 6
$ws.Range("B34").Value = "audio"

$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "audio"

$ws.Range("A36").Value = 7
$ws.Range("B36").Value = "audio"
